# #5: cash & deposit done
#
# The "存款" (deposit) worksheet is extended from a bare A1:F11 table into
# the same A1:M<n> schema used by the other property-type worksheets
# (property_category / category / date / legislator_name / legislator_id /
# source_file / index columns appended), and its header row (row 1) is
# turned into real column-name labels instead of being a duplicate of the
# first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 存款

# ---- header row -----------------------------------------------------
$headers = @{
    "B1" = "bank";
    "C1" = "deposit_type";
    "D1" = "currency";
    "E1" = "owner";
    "F1" = "total";
    "G1" = "property_category";
    "H1" = "category";
    "I1" = "date";
    "J1" = "legislator_name";
    "K1" = "legislator_id";
    "L1" = "source_file";
    "M1" = "index";
}

foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value2 = $headers[$addr]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(1).LineStyle = 1
    $cell.Borders.Item(2).LineStyle = 1
    $cell.Borders.Item(3).LineStyle = 1
    $cell.Borders.Item(4).LineStyle = 1
}

# ---- data rows --------------------------------------------------------
$rows = @(
    @{ r=2;  A=60; B="陽信商業銀行民生分行";           C="定期存款";    D="新臺幣"; E="饒月琴"; F=2200000 },
    @{ r=3;  A=61; B="國泰世華商業銀行三民分行";         C="定期存款";    D="新臺幣"; E="饒月琴"; F=3100000 },
    @{ r=4;  A=62; B="國泰世華商業銀行三民分行";         C="活期儲蓄存款"; D="新臺幣"; E="饒月琴"; F=6751 },
    @{ r=5;  A=63; B="中華郵政股份有限公司台北民生郵局"; C="活期存款";    D="新臺幣"; E="饒月琴"; F=12426 },
    @{ r=6;  A=64; B="台新國際商業銀行";                 C="綜合存款";    D="新臺幣"; E="饒月琴"; F=111 },
    @{ r=7;  A=65; B="陽信商業銀行民生分行";             C="活期儲蓄存款"; D="新臺幣"; E="饒月琴"; F=0 },
    @{ r=8;  A=66; B="中華郵政股份有限公司";             C="活期存款";    D="新臺幣"; E="許忠信"; F=829818 },
    @{ r=9;  A=67; B="中華郵政股份有限公司";             C="定期存款";    D="新臺幣"; E="許忠信"; F=343559 },
    @{ r=10; A=68; B="兆豐國際商業銀行";                 C="綜合存款";    D="新臺幣"; E="許忠信"; F=11791 },
    @{ r=11; A=69; B="國泰世華商業銀行三民分行";         C="活期儲蓄存款"; D="新臺幣"; E="許忠信"; F=116996 }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Range("A$r").Value2 = $row.A
    $ws.Range("B$r").Value2 = $row.B
    $ws.Range("C$r").Value2 = $row.C
    $ws.Range("D$r").Value2 = $row.D
    $ws.Range("E$r").Value2 = $row.E
    $ws.Range("F$r").Value2 = $row.F

    # newly appended constant columns shared with every other property sheet
    $ws.Range("G$r").Value2 = "deposit"
    $ws.Range("H$r").Value2 = "normal"
    $ws.Range("I$r").Value2 = "2012-04-23"
    $ws.Range("J$r").Value2 = "許忠信"
    $ws.Range("K$r").Value2 = 1749
    $ws.Range("L$r").Value2 = "tmp50641"
    $ws.Range("M$r").Value2 = $row.A
}
